$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 6 de Octubre de 2020 a las 03:00"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 7678459
$ws.Range("C4").Value = 40391
$ws.Range("D4").Value = 4890306
$ws.Range("E4").Value = 2573159
$ws.Range("G4").Value = 383
$ws.Range("H4").Value = 214994

# Row 11 - Argentina
$ws.Range("B11").Value = 809728
$ws.Range("C11").Value = 11242
$ws.Range("D11").Value = 649017
$ws.Range("E11").Value = 139243
$ws.Range("G11").Value = 450
$ws.Range("H11").Value = 21468

# Row 29 - Canada
$ws.Range("B29").Value = 168960
$ws.Range("C29").Value = 2804
$ws.Range("D29").Value = 142334
$ws.Range("E29").Value = 17122
$ws.Range("G29").Value = 23
$ws.Range("H29").Value = 9504

# Row 69 - Paraguay
$ws.Range("B69").Value = 44715
$ws.Range("C69").Value = 533
$ws.Range("D69").Value = 27887
$ws.Range("E69").Value = 15881
$ws.Range("G69").Value = 18
$ws.Range("H69").Value = 947

# Row 192 - Bermudas
$ws.Range("D192").Value = 170
$ws.Range("E192").Value = 2
